$d = $word.ActiveDocument

# Locate the merged line: "237 - 7 Eidolons, Incendiary Ichor 432 - 4 Bloat Flies, Unhinging Jaw "
# It needs to become two separate paragraphs, split at the single space
# that currently joins "Ichor" and "432" (that space is removed, replaced
# by a paragraph break) so totals aren't double-counted on one line.
$anchor = "Incendiary Ichor 432 - 4 Bloat Flies, Unhinging Jaw"
$found = $d.Content
$found.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# The space sits right after "Incendiary Ichor" inside the matched range.
$splitAt = $found.Start + "Incendiary Ichor".Length
$spaceRng = $d.Range($splitAt, $splitAt + 1)

# Remove the joining space, then turn that point into a new paragraph break.
$spaceRng.Text = ""
$spaceRng.Collapse(1)
$spaceRng.InsertParagraphBefore()
